# Converts R,G,B (0-255) into the little-endian OLE COLORREF integer that
# Excel's Interior.Color property expects.
function Get-OleColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$wb = $excel.ActiveWorkbook

# Duplicate the "WebTable" worksheet. Excel inserts the copy immediately
# after the source sheet, i.e. right after "WebTable".
$webTable = $wb.Worksheets.Item("WebTable")
$webTable.Copy($null, $webTable)

$newSheet = $wb.Worksheets.Item(3)
$newSheet.Name = "newSheet"

# Update the header row text on the new sheet.
$newSheet.Range("B1").Value = "Contact"
$newSheet.Range("C1").Value = "Country"

# Shade the header row (A1:C1), one cell at a time, with the same light
# indigo/gray tone (standard palette color index 55).
$newSheet.Range("A1").Interior.ColorIndex = 55
$newSheet.Range("B1").Interior.Color = (Get-OleColor 0x34 0x33 0x99)
$newSheet.Range("C1").Interior.Color = (Get-OleColor 0x35 0x33 0x99)

# The new sheet's selection collapses to a single cell, A2.
$newSheet.Range("A2").Select()

# "newSheet" becomes the active tab (the 3rd sheet, 0-based activeTab = 2).
$newSheet.Activate()
